$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: extend bottom-border-only formatting into new column L (empty cell)
$ws.Range("K3").Copy()
$ws.Range("L3").PasteSpecial(-4122)

# Row 4: new year column (2021), matching the style of the existing year cells
$ws.Range("K4").Copy()
$ws.Range("L4").PasteSpecial(-4122)
$ws.Range("L4").Value = 2021

# Row 5: new data value (269), matching the style of the existing data cells
$ws.Range("K5").Copy()
$ws.Range("L5").PasteSpecial(-4122)
$ws.Range("L5").Value = 269

# Update the active selection to reflect where the editor left off
$ws.Range("N3").Select() | Out-Null
